# Regional Availability Factor.xlsx - "updated 4.0 files and mdl"
#
# 1. About sheet: bump the "last updated" date from 3/15/2024 to 3/28/2024
# 2. RAF-capacity sheet: raise the RAF for the two hydrogen technologies
#    (hydrogen combustion turbine, hydrogen combined cycle) from 0.3 to 1
# 3. Leave the workbook with the RAF-capacity sheet active/selected, scrolled
#    down and zoomed to 80%, matching the author's final view state.

$wb = $excel.ActiveWorkbook

# --- About sheet: update date in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = Get-Date -Year 2024 -Month 3 -Day 28 -Hour 0 -Minute 0 -Second 0

# --- RAF-capacity sheet: update hydrogen RAF values ---
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# --- Make RAF-capacity the active sheet / set its view state ---
$wsCapacity.Activate()
$excel.Goto($wsCapacity.Range("A14"), $true)
$wsCapacity.Range("B25").Select()
$excel.ActiveWindow.Zoom = 80
